$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2352347881312385
$ws.Range("C2").Value = 0.7011352911275793
$ws.Range("D2").Value = 0.754111145231435
$ws.Range("E2").Value = 0.868395730776836
$ws.Range("F2").Value = 0.8601628886204244
$ws.Range("G2").Value = 18

$ws.Range("B3").Value = 0.2672786791899959
$ws.Range("C3").Value = 0.45307994943011
$ws.Range("D3").Value = 0.2955295970660748
$ws.Range("E3").Value = 0.5436263395624561
$ws.Range("F3").Value = 0.4879522889190182
$ws.Range("G3").Value = 17

$ws.Range("B4").Value = 0.3111056731803871
$ws.Range("C4").Value = 0.3624139282624727
$ws.Range("D4").Value = 0.1723564222806536
$ws.Range("E4").Value = 0.4151583099019621
$ws.Range("F4").Value = 0.2839148837627348
$ws.Range("G4").Value = 16

$ws.Range("B5").Value = 0.3455294579606332
$ws.Range("C5").Value = 0.3658231501549597
$ws.Range("D5").Value = 0.1878025140189348
$ws.Range("E5").Value = 0.4333618742101511
$ws.Range("F5").Value = 0.2707369064906118
$ws.Range("G5").Value = 15

$ws.Range("B6").Value = 0.3584323639489301
$ws.Range("C6").Value = 0.3839945786029612
$ws.Range("D6").Value = 0.2002730478501564
$ws.Range("E6").Value = 0.4475187681540925
$ws.Range("F6").Value = 0.2780688952452576
$ws.Range("G6").Value = 14

$ws.Range("B7").Value = 0.364459496889083
$ws.Range("C7").Value = 0.4106453173185426
$ws.Range("D7").Value = 0.2107924313882408
$ws.Range("E7").Value = 0.4591213689083103
$ws.Range("F7").Value = 0.2906174726197602
$ws.Range("G7").Value = 13

$ws.Range("B8").Value = 0.4143658927697875
$ws.Range("C8").Value = 0.4310987738402203
$ws.Range("D8").Value = 0.221430687247661
$ws.Range("E8").Value = 0.4705642222350324
$ws.Range("F8").Value = 0.2329219787203619
$ws.Range("G8").Value = 12

$ws.Range("B9").Value = 0.3869142407525709
$ws.Range("C9").Value = 0.4001865847611734
$ws.Range("D9").Value = 0.1924388048854656
$ws.Range("E9").Value = 0.4386784755210422
$ws.Range("F9").Value = 0.2168174179053887
$ws.Range("G9").Value = 11

$ws.Range("B10").Value = 0.3617893501639269
$ws.Range("C10").Value = 0.3734161284103805
$ws.Range("D10").Value = 0.1744359538854561
$ws.Range("E10").Value = 0.4176553051087178
$ws.Range("F10").Value = 0.2199606530304394
$ws.Range("G10").Value = 10

$ws.Range("B11").Value = 0.3719529320823554
$ws.Range("C11").Value = 0.3812107115717386
$ws.Range("D11").Value = 0.1830055333568115
$ws.Range("E11").Value = 0.4277914601260893
$ws.Range("F11").Value = 0.2241397295910946
$ws.Range("G11").Value = 9

